$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint backlog restructuring ---
# Column C task labels move into column D; column A keeps a short " Name"
# label for grouped rows, while the task text now lives once in column D.

# "Identifie 3 code smells" block (rows 4-8)
$ws.Cells.Item(4,3).ClearContents()
$ws.Cells.Item(4,4).Value = "Identifie 3 code smells"

$ws.Cells.Item(5,1).Value = " Ana Gadelha"
$ws.Cells.Item(5,4).Value = "Identifie 3 code smells"

$ws.Cells.Item(6,1).Value = " Rodrigo Mesquita"
$ws.Cells.Item(6,4).Value = "Identifie 3 code smells"

$ws.Cells.Item(7,3).ClearContents()
$ws.Cells.Item(7,4).Value = "Identifie 3 code smells"

$ws.Cells.Item(8,3).ClearContents()
$ws.Cells.Item(8,4).Value = "Identifie 3 code smells"

# "Identifie 3 design paterns" block (rows 10-14)
$ws.Cells.Item(10,3).ClearContents()
$ws.Cells.Item(10,4).Value = "Identifie 3 design paterns"

$ws.Cells.Item(11,1).Value = " Ana Gadelha"
$ws.Cells.Item(11,4).Value = "Identifie 3 design paterns"

$ws.Cells.Item(12,1).Value = " Rodrigo Mesquita"
$ws.Cells.Item(12,4).Value = "Identifie 3 design paterns"

$ws.Cells.Item(13,3).ClearContents()
$ws.Cells.Item(13,4).Value = "Identifie 3 design paterns"

$ws.Cells.Item(14,3).ClearContents()
$ws.Cells.Item(14,4).Value = "Identifie 3 design paterns "

# "Review each others design patterns - All" row (row 16)
$ws.Cells.Item(16,1).Value = " All"
$ws.Cells.Item(16,4).Value = "Review each others design patterns "

# --- Column width / selection cosmetic changes ---
# (ColumnWidth is stored in 1/6-character increments by the host, so 31
#  is the closest achievable value to the authored 31.88671875 width)
$ws.Columns.Item(4).ColumnWidth = 31
$ws.Range("E18").Select() | Out-Null
